$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETS")

# Electricity capacity additions calibration:
# bump the "natural gas nonpeaker" and "onshore wind" shareweights across all
# forecast years (2020-2050, columns B:AF).
$ws.Range("B3:AF3").Value = 10   # natural gas nonpeaker: 1 -> 10
$ws.Range("B6:AF6").Value = 3    # onshore wind: 1 -> 3

# Leave the selection where the edit ended up, matching the saved view.
$ws.Range("B6:AF6").Select()

# Restore "About" as the active sheet/tab (it was active before the edit).
$wb.Worksheets.Item("About").Activate()
